$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 107, shifting existing rows 107..210 down to 108..211
# (the sheet's dimension grows from A1:T210 to A1:T211).
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row 107 with the new price-record data.
$ws.Range("A107").Value() = 4
$ws.Range("B107").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C107").Value() = "Los Lagos"
$ws.Range("D107").Value() = 44589
$ws.Range("E107").Value() = 10
$ws.Range("F107").Value() = "Fruta"
$ws.Range("G107").Value() = 100104
$ws.Range("H107").Value() = "Frutos de pepita"
$ws.Range("I107").Value() = 100104005
$ws.Range("J107").Value() = "Pera"
$ws.Range("K107").Value() = "Packham's Triumph"
$ws.Range("L107").Value() = "Primera"
$ws.Range("M107").Value() = 500
$ws.Range("N107").Value() = 14000
$ws.Range("O107").Value() = 15000
$ws.Range("P107").Value() = 14500
$ws.Range("Q107").Value() = "$/caja 15 kilos empedrada"
$ws.Range("R107").Value() = "Región de O'Higgins"
$ws.Range("S107").Value() = 967
$ws.Range("T107").Value() = 15
